# Uprava vzorce pro vypocet stravneho za sluzebni cesty
# Remove the "Km" column contents (header + per-row values + the SUM totals row)
# from both worksheets, then leave the selection/active sheet state matching
# the author's final interaction (List2 active, specific cell selections).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- List1 ---
$ws1.Range("G1").ClearContents() | Out-Null
$ws1.Range("G10").ClearContents() | Out-Null
$ws1.Range("G19").ClearContents() | Out-Null
$ws1.Range("G21").ClearContents() | Out-Null
$ws1.Rows.Item(32).Delete() | Out-Null

# --- List2 ---
$ws2.Range("G1").ClearContents() | Out-Null
$ws2.Range("G5").ClearContents() | Out-Null
$ws2.Range("G19").ClearContents() | Out-Null
$ws2.Range("G20").ClearContents() | Out-Null
$ws2.Range("G21").ClearContents() | Out-Null
$ws2.Rows.Item(32).Delete() | Out-Null

# --- Final selection / active sheet state ---
$ws1.Range("J8").Select() | Out-Null
$ws2.Range("H11").Select() | Out-Null
